$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert 4 new rows (a new date's worth of quality grades)
# above the existing data block that starts at row 250, pushing the old
# rows 250-256 down to 254-260.
$ws.Rows("250:253").Insert()

# Columns that are constant across every record in this block.
$ws.Range("A250:A253").Value2 = 2
$ws.Range("B250:B253").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C250:C253").Value2 = "Coquimbo"
$ws.Range("D250:D253").Value2 = 44461
$ws.Range("E250:E253").Value2 = 4
$ws.Range("F250:F253").Value2 = 100112043
$ws.Range("G250:G253").Value2 = "Pepino dulce"
$ws.Range("H250:H253").Value2 = "Cultivar IV Región"
$ws.Range("N250:N253").Value2 = "`$/bandeja 18 kilos"
$ws.Range("O250:O253").Value2 = "Provincia de Limarí"
$ws.Range("Q250:Q253").Value2 = 18
$ws.Range("R250:R253").Value2 = "Hortaliza"

# Row 250: Especial
$ws.Range("I250").Value2 = "Especial"
$ws.Range("J250").Value2 = 300
$ws.Range("K250").Value2 = 13500
$ws.Range("L250").Value2 = 14000
$ws.Range("M250").Value2 = 13750
$ws.Range("P250").Value2 = 764

# Row 251: Primera
$ws.Range("I251").Value2 = "Primera"
$ws.Range("J251").Value2 = 400
$ws.Range("K251").Value2 = 11500
$ws.Range("L251").Value2 = 12000
$ws.Range("M251").Value2 = 11750
$ws.Range("P251").Value2 = 653

# Row 252: Segunda
$ws.Range("I252").Value2 = "Segunda"
$ws.Range("J252").Value2 = 400
$ws.Range("K252").Value2 = 9500
$ws.Range("L252").Value2 = 10000
$ws.Range("M252").Value2 = 9750
$ws.Range("P252").Value2 = 542

# Row 253: Tercera
$ws.Range("I253").Value2 = "Tercera"
$ws.Range("J253").Value2 = 300
$ws.Range("K253").Value2 = 6500
$ws.Range("L253").Value2 = 7000
$ws.Range("M253").Value2 = 6750
$ws.Range("P253").Value2 = 375
